$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NroPoliza value (E2): 12112002068 -> 12112002080
$ws.Range("E2").Value = 12112002080

# Update FechaSiniestro value (G2): 23/06/2022 -> 19/03/2021
# Leading apostrophe preserves the existing "stored as text" (quotePrefix) formatting
$ws.Range("G2").Value = "'19/03/2021"

# Remove only the hyperlink attached to B2, leaving the other hyperlinks untouched
$target = $ws.Range("B2")
$toDelete = @()
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Row -eq $target.Row -and $h.Range.Column -eq $target.Column) {
        $toDelete += $h
    }
}
foreach ($d in $toDelete) {
    $d.Delete()
}

# Update the active selection / view
$ws.Range("G3").Select()
